$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3: ReportPath -> ReportPath_Bot1, value -> Bot 1 subfolder
$ws.Range("B3").Value = "D:\ReportsScraping\Bot 1\"
$ws.Range("A3").Value = "ReportPath_Bot1"

# Add new row 4: ReportPath_main
$ws.Range("A4").Value = "ReportPath_main"
$ws.Range("B4").Value = "D:\ReportsScraping\MAIN.xlsx"

# Update selection to G3
$ws.Range("G3").Select()
